# Update odds values for existing matches (rows 2,3,5,6,7,8,9,10,11)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 17
$ws.Range("Q2").Value = 1.53
$ws.Range("R2").Value = 2.5
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = 1.36
$ws.Range("Q3").Value = 2.1
$ws.Range("R3").Value = 1.73
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.91
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 4
$ws.Range("L5").Value = 4.75
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 2.4
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.67
$ws.Range("W5").Value = 6
$ws.Range("X5").Value = 8.5
$ws.Range("AC5").Value = 7
$ws.Range("AE5").Value = 19
$ws.Range("AF5").Value = 67
$ws.Range("AI5").Value = 19
$ws.Range("AK5").Value = 41
$ws.Range("AO5").Value = 12
$ws.Range("AT5").Value = 2.38
$ws.Range("AU5").Value = 9
$ws.Range("AV5").Value = 67
$ws.Range("AX5").Value = 5.5
$ws.Range("AY5").Value = 23
$ws.Range("BA5").Value = 81
$ws.Range("BB5").Value = 126
$ws.Range("G6").Value = 1.82
$ws.Range("G7").Value = 1.37
$ws.Range("Q7").Value = 1.44
$ws.Range("U7").Value = 1.57
$ws.Range("H8").Value = 3.6
$ws.Range("J8").Value = 3
$ws.Range("Q8").Value = 1.53
$ws.Range("R8").Value = 2.4
$ws.Range("S8").Value = 1.29
$ws.Range("T8").Value = 3.5
$ws.Range("U8").Value = 1.44
$ws.Range("V8").Value = 2.63
$ws.Range("AG8").Value = 101
$ws.Range("AH8").Value = 13
$ws.Range("AJ8").Value = 10
$ws.Range("AK8").Value = 26
$ws.Range("AO8").Value = 13
$ws.Range("AP8").Value = 19
$ws.Range("AS8").Value = 101
$ws.Range("AT8").Value = 3.5
$ws.Range("G9").Value = 1.33
$ws.Range("H9").Value = 6.25
$ws.Range("I9").Value = 6.25
$ws.Range("J9").Value = 1.73
$ws.Range("L9").Value = 5.5
$ws.Range("O9").Value = 1.06
$ws.Range("P9").Value = 10
$ws.Range("Q9").Value = 1.22
$ws.Range("U9").Value = 1.4
$ws.Range("V9").Value = 2.75
$ws.Range("Z9").Value = 13
$ws.Range("AB9").Value = 15
$ws.Range("AI9").Value = 41
$ws.Range("AJ9").Value = 21
$ws.Range("AM9").Value = 29
$ws.Range("AO9").Value = 6.5
$ws.Range("AZ9").Value = 23
$ws.Range("BB9").Value = 67
$ws.Range("Q10").Value = 1.5
$ws.Range("U10").Value = 1.53
$ws.Range("V10").Value = 2.38
$ws.Range("G11").Value = 2.75
$ws.Range("I11").Value = 2.25
$ws.Range("J11").Value = 3.1
$ws.Range("K11").Value = 2.6
$ws.Range("N11").Value = 26
$ws.Range("Q11").Value = 1.36
$ws.Range("R11").Value = 3.1
$ws.Range("U11").Value = 1.33
$ws.Range("AA11").Value = 19
$ws.Range("AC11").Value = 26
$ws.Range("AE11").Value = 11
$ws.Range("AF11").Value = 26
$ws.Range("AI11").Value = 17
$ws.Range("AJ11").Value = 10
$ws.Range("AK11").Value = 23
$ws.Range("AP11").Value = 17
$ws.Range("AS11").Value = 81
$ws.Range("AY11").Value = 11
# Remove the two matches that dropped out of this week's list:
#   row 12: Puszcza vs Widzew Lodz (POLAND - EKSTRAKLASA)
#   row 13: Legnica vs Wisla Plock (POLAND - DIVISION 1)
# Deleting shifts row 14 (UTA Arad vs Univ. Craiova) up to become the new row 12.
$ws.Rows(12).EntireRow.Delete()
$ws.Rows(12).EntireRow.Delete()

# A handful of odds on the (now-shifted) UTA Arad vs Univ. Craiova row were
# also updated as part of this refresh.
$ws.Range("I12").Value = 1.87
$ws.Range("O12").Value = 1.29
$ws.Range("P12").Value = 3.5
$ws.Range("Q12").Value = 1.98
$ws.Range("R12").Value = 1.88
